$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.191709995269775
$ws.Range("B1").Value = 2.372738122940063
$ws.Range("C1").Value = 4.246733665466309
$ws.Range("D1").Value = 2.87806224822998
$ws.Range("E1").Value = 1.120337724685669
